$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 97 ---
# Copy formatting from the row above (row 96) first so the new row's
# cells inherit the same styles (e.g. the date/time number format on
# column A) without Excel creating brand-new style records.
$ws.Range("A96:H96").Copy($ws.Range("A97:H97"))

$ws.Cells.Item(97, 1).Value = 45454.2916666667
$ws.Cells.Item(97, 2).Value = 0
$ws.Cells.Item(97, 3).Value = 2
$ws.Cells.Item(97, 4).Value = 2
$ws.Cells.Item(97, 5).Value = 2
$ws.Cells.Item(97, 6).Value = 2

# Column G ("adj_close") holds the text "2" (not the number 2) in the
# source data, so force text storage via the "@" number format, then
# restore the default style so no stray formatting is left behind.
$ws.Cells.Item(97, 7).NumberFormat = "@"
$ws.Cells.Item(97, 7).Value = "2"
$ws.Cells.Item(97, 7).Style = "Normal"

$ws.Cells.Item(97, 8).Value = "KK.MI"

# --- Row 98 ---
$ws.Range("A96:H96").Copy($ws.Range("A98:H98"))

$ws.Cells.Item(98, 1).Value = 45455.2916666667
$ws.Cells.Item(98, 2).Value = 0
$ws.Cells.Item(98, 3).Value = 2
$ws.Cells.Item(98, 4).Value = 2
$ws.Cells.Item(98, 5).Value = 2
$ws.Cells.Item(98, 6).Value = 2

$ws.Cells.Item(98, 7).NumberFormat = "@"
$ws.Cells.Item(98, 7).Value = "2"
$ws.Cells.Item(98, 7).Style = "Normal"

$ws.Cells.Item(98, 8).Value = "KK.MI"
